# "Se agrega talle a producto" - mark the "agregar talle a prod" task as
# done (100%), push the "en proceso" status down to the next task
# (zapatilla por talle), and flag "cambiar precio de producto" as in
# progress too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 23: "agregar talle a prod" -> completed (100%)
$ws.Range("C23").Value = 1
$ws.Range("C23").NumberFormat = "0%"
$ws.Range("C23").Font.Underline = $false

# Row 24: "reporte de zapatilla por talle en stock..." -> now in progress
$ws.Range("C24").Value = "en proceso"

# Row 37: "cambiar precio de producto" -> now in progress (no longer underlined)
$ws.Range("C37").Value = "en proceso"
$ws.Range("C37").Font.Underline = $false

# Update the view state to match (scrolled/selected a couple rows down)
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("C25").Select()
